$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D from 28 to 30 (raw OOXML width units). The COM ColumnWidth
# property applies Excel's internal padding offset (~0.8333) on top of the
# value you set, so back that out to land on an exact stored width of 30.
$ws.Columns.Item(4).ColumnWidth = 29.16666666666667

# Insert a new row above the current row 5; this shifts rows 5-8 down to
# 6-9 (values, number formats and the hyperlink style all move with them).
$ws.Rows.Item(5).Insert()

# New row 5 content (a freshly scraped listing).
$ws.Range("A5").Value = "2026-01-10 12:47:39"
$ws.Range("B5").Value = "【急募】WebシステムのJavaScript/PHP開発エンジニア募集"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5468735"
$ws.Range("G5").Value = 185
$ws.Range("H5").Value = "★Java ◆開発 ○PHP"

# The whole sheet was re-scraped at 12:47:39, so every row's "retrieved at"
# timestamp in column A advances to the new run time (not just the rows
# that physically shifted down).
$ws.Range("A2").Value = "2026-01-10 12:47:39"
$ws.Range("A3").Value = "2026-01-10 12:47:39"
$ws.Range("A4").Value = "2026-01-10 12:47:39"
$ws.Range("A6").Value = "2026-01-10 12:47:39"
$ws.Range("A7").Value = "2026-01-10 12:47:39"
$ws.Range("A8").Value = "2026-01-10 12:47:39"
$ws.Range("A9").Value = "2026-01-10 12:47:39"

# The row-insert shifted cell contents but left the worksheet Hyperlinks
# collection's ref->target mapping out of sync with the new row positions,
# so rebuild it from scratch in row order.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5468493")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5468303")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5468735")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5468677")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5468743")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5468432")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5468565")
